# Applies the Tab14 annex-table update:
#  - "South Sudan" -> "South Sudan*" (and gains the footnote highlight fill)
#  - "Nigeria*" -> "Nigeria" (and loses the footnote highlight fill)
#  - refreshed indicator figures for a block of country rows (69-98)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab14")

# --- Country name / footnote-marker changes -------------------------------

$ws.Range("B34").Value2 = "South Sudan*"
$ws.Range("B57").Value2 = "Nigeria"

# South Sudan (row 34) now carries the "*" footnote, so it picks up the same
# light-blue highlight fill used by the other starred countries (e.g. row 17,
# Chad*).
$ws.Range("B17:H17").Copy()
$ws.Range("B34:H34").PasteSpecial(-4122)  # xlPasteFormats

# Nigeria (row 57) no longer carries the "*" footnote, so it reverts to the
# plain (unhighlighted) row formatting used by the rest of the table (e.g.
# row 5, Eswatini).
$ws.Range("B5:H5").Copy()
$ws.Range("B57:H57").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Refreshed data values --------------------------------------------------

$ws.Range("C69").Value2 = 3.87775271279471
$ws.Range("D69").Value2 = 0.38498073390552
$ws.Range("E69").Value2 = 0.66522204875946
$ws.Range("F69").Value2 = 0.69935835685049
$ws.Range("G69").Value2 = 0.13515450463941
$ws.Range("H69").Value2 = 0.63066862736429

$ws.Range("C77").Value2 = 5.86061244540745
$ws.Range("D77").Value2 = 0.34672289755609
$ws.Range("E77").Value2 = 0.79228505823347
$ws.Range("F77").Value2 = 0.8372622595893
$ws.Range("G77").Value2 = -0.0692543857214
$ws.Range("H77").Value2 = 0.84942626290851

$ws.Range("C80").Value2 = 4.4020930826664
$ws.Range("D80").Value2 = 0.37987058795989
$ws.Range("E80").Value2 = 0.63435044139624
$ws.Range("F80").Value2 = 0.57064414024353
$ws.Range("G80").Value2 = -0.0758193990083
$ws.Range("H80").Value2 = 0.72097374498844

$ws.Range("C82").Value2 = 4.52134836025727
$ws.Range("D82").Value2 = 0.34189509810546
$ws.Range("E82").Value2 = 0.68167668886674
$ws.Range("F82").Value2 = 0.71937807935935
$ws.Range("G82").Value2 = 0.02317824886229
$ws.Range("H82").Value2 = 0.68049986163775

$ws.Range("C84").Value2 = 4.33113220604983
$ws.Range("D84").Value2 = 0.38575830716978
$ws.Range("E84").Value2 = 0.65880039063367
$ws.Range("F84").Value2 = 0.68839752945033
$ws.Range("G84").Value2 = 0.06544675915502
$ws.Range("H84").Value2 = 0.64880306205966

$ws.Range("C86").Value2 = 4.58608476739181
$ws.Range("D86").Value2 = 0.3220815917379
$ws.Range("E86").Value2 = 0.67463707296472
$ws.Range("F86").Value2 = 0.67789838972845
$ws.Range("G86").Value2 = -0.0152199844045
$ws.Range("H86").Value2 = 0.69189433361355

$ws.Range("C87").Value2 = 5.13766500353813
$ws.Range("D87").Value2 = 0.33295808297892
$ws.Range("E87").Value2 = 0.73048964763681
$ws.Range("F87").Value2 = 0.81976922725638
$ws.Range("G87").Value2 = 0.11536251547902
$ws.Range("H87").Value2 = 0.78173424055179

$ws.Range("C89").Value2 = 5.61383240150683
$ws.Range("D89").Value2 = 0.30502390954643
$ws.Range("E89").Value2 = 0.69966786634177
$ws.Range("F89").Value2 = 0.81702966581691
$ws.Range("G89").Value2 = -0.0383586944222
$ws.Range("H89").Value2 = 0.83823045636668

$ws.Range("C90").Value2 = 6.659212203736
$ws.Range("D90").Value2 = 0.25483741805605
$ws.Range("E90").Value2 = 0.74898661608281
$ws.Range("F90").Value2 = 0.86344047473825
$ws.Range("G90").Value2 = -0.0376309784959
$ws.Range("H90").Value2 = 0.91528327698293

$ws.Range("C94").Value2 = 5.54514534132821
$ws.Range("D94").Value2 = 0.25595080852509
$ws.Range("E94").Value2 = 0.73234396179517
$ws.Range("F94").Value2 = 0.83461324657713
$ws.Range("G94").Value2 = 0.03499539941549
$ws.Range("H94").Value2 = 0.81901471103941

$ws.Range("C97").Value2 = 4.4427880118875
$ws.Range("D97").Value2 = 0.36188379850458
$ws.Range("E97").Value2 = 0.67324354894021
$ws.Range("F97").Value2 = 0.67438138495473
$ws.Range("G97").Value2 = 0.03863246952824
$ws.Range("H97").Value2 = 0.68029793571023

$ws.Range("C98").Value2 = 4.75531077384949
$ws.Range("D98").Value2 = 0.37425981724963
$ws.Range("E98").Value2 = 0.66255626433036
$ws.Range("F98").Value2 = 0.73614815228126
$ws.Range("G98").Value2 = 0.07607171985814
$ws.Range("H98").Value2 = 0.72445550736259
